$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value updates (refreshed crypto price/volume snapshot).
# NumberFormat "@" (Text) forces Excel to keep the literal string instead of
# re-parsing numeric-looking text (e.g. "1.0000" -> 1); ClearFormats() afterwards
# drops the temporary style so cell formatting stays exactly as it was before.
$updates = @(
    @{ Cell = 'D2'; Value = '23.793.43' }
    @{ Cell = 'E2'; Value = '  -0.89%  ' }
    @{ Cell = 'D3'; Value = '1.630.36' }
    @{ Cell = 'E3'; Value = '  -1.66%  ' }
    @{ Cell = 'D4'; Value = '1.002' }
    @{ Cell = 'E4'; Value = '  -0.06%  ' }
    @{ Cell = 'D5'; Value = '307.82' }
    @{ Cell = 'E5'; Value = '  -0.78%  ' }
    @{ Cell = 'D6'; Value = '1.002' }
    @{ Cell = 'E6'; Value = '  -0.01%  ' }
    @{ Cell = 'D7'; Value = '0.3829' }
    @{ Cell = 'E7'; Value = '  -1.84%  ' }
    @{ Cell = 'D8'; Value = '0.3803' }
    @{ Cell = 'E8'; Value = '  -1.76%  ' }
    @{ Cell = 'D9'; Value = '50.50' }
    @{ Cell = 'E9'; Value = '  -1.23%  ' }
    @{ Cell = 'D10'; Value = '1.311' }
    @{ Cell = 'E10'; Value = '  -4.16%  ' }
    @{ Cell = 'D11'; Value = '1.002' }
    @{ Cell = 'E11'; Value = '  -0.02%  ' }
    @{ Cell = 'D12'; Value = '0.08319' }
    @{ Cell = 'E12'; Value = '  -2.08%  ' }
    @{ Cell = 'D13'; Value = '23.53' }
    @{ Cell = 'E13'; Value = '  -1.96%  ' }
    @{ Cell = 'D14'; Value = '6.904' }
    @{ Cell = 'E14'; Value = '  -4.28%  ' }
    @{ Cell = 'D15'; Value = '7.698' }
    @{ Cell = 'E15'; Value = '  -3.99%  ' }
    @{ Cell = 'D16'; Value = '0.00001287' }
    @{ Cell = 'E16'; Value = '  -2.16%  ' }
    @{ Cell = 'D17'; Value = '1.635.86' }
    @{ Cell = 'E17'; Value = '  -1.08%  ' }
    @{ Cell = 'D18'; Value = '93.09' }
    @{ Cell = 'E18'; Value = '  -1.61%  ' }
    @{ Cell = 'D19'; Value = '0.06919' }
    @{ Cell = 'E19'; Value = '  -1.03%  ' }
    @{ Cell = 'D20'; Value = '19.26' }
    @{ Cell = 'E20'; Value = '  -3.69%  ' }
    @{ Cell = 'D21'; Value = '6.807' }
    @{ Cell = 'E21'; Value = '  -2.71%  ' }
    @{ Cell = 'D22'; Value = '1.0000' }
    @{ Cell = 'E22'; Value = '  -0.21%  ' }
    @{ Cell = 'D23'; Value = '13.43' }
    @{ Cell = 'E23'; Value = '  -2.06%  ' }
    @{ Cell = 'D24'; Value = '23.788.14' }
    @{ Cell = 'E24'; Value = '  -0.92%  ' }
    @{ Cell = 'D25'; Value = '2.429' }
    @{ Cell = 'E25'; Value = '  -2.34%  ' }
    @{ Cell = 'D26'; Value = '2.834' }
    @{ Cell = 'E26'; Value = '  -9.29%  ' }
    @{ Cell = 'D27'; Value = '21.72' }
    @{ Cell = 'E27'; Value = '  -2.54%  ' }
    @{ Cell = 'D28'; Value = '152.14' }
    @{ Cell = 'E28'; Value = '  -1.00%  ' }
    @{ Cell = 'D29'; Value = '5.449' }
    @{ Cell = 'E29'; Value = '  +2.42%  ' }
    @{ Cell = 'D30'; Value = '135.55' }
    @{ Cell = 'E30'; Value = '  -3.64%  ' }
    @{ Cell = 'D31'; Value = '7.840' }
    @{ Cell = 'E31'; Value = '  +0.24%  ' }
    @{ Cell = 'D32'; Value = '2.479' }
    @{ Cell = 'E32'; Value = '  -0.69%  ' }
    @{ Cell = 'D33'; Value = '1.816.55' }
    @{ Cell = 'E33'; Value = '  -0.95%  ' }
    @{ Cell = 'B34'; Value = 'ImmutableX' }
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' }
    @{ Cell = 'D34'; Value = '0.9753' }
    @{ Cell = 'E34'; Value = '  -6.99%  ' }
    @{ Cell = 'B35'; Value = 'Hedera' }
    @{ Cell = 'C35'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' }
    @{ Cell = 'D35'; Value = '0.07876' }
    @{ Cell = 'E35'; Value = '  -3.53%  ' }
    @{ Cell = 'D36'; Value = '0.02863' }
    @{ Cell = 'E36'; Value = '  -4.86%  ' }
    @{ Cell = 'D37'; Value = '6.527' }
    @{ Cell = 'E37'; Value = '  -2.35%  ' }
    @{ Cell = 'D38'; Value = '0.2634' }
    @{ Cell = 'E38'; Value = '  -2.73%  ' }
    @{ Cell = 'D39'; Value = '10.33' }
    @{ Cell = 'E39'; Value = '  -7.84%  ' }
    @{ Cell = 'D40'; Value = '0.09020' }
    @{ Cell = 'E40'; Value = '  -1.41%  ' }
    @{ Cell = 'D41'; Value = '0.7410' }
    @{ Cell = 'E41'; Value = '  -2.28%  ' }
    @{ Cell = 'B42'; Value = 'TrustWalletToken' }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt' }
    @{ Cell = 'D42'; Value = '1.410' }
    @{ Cell = 'E42'; Value = '  -0.91%  ' }
    @{ Cell = 'B43'; Value = 'Aptos' }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt' }
    @{ Cell = 'D43'; Value = '13.10' }
    @{ Cell = 'E43'; Value = '  -3.86%  ' }
    @{ Cell = 'D44'; Value = '16.35' }
    @{ Cell = 'E44'; Value = '  -1.55%  ' }
    @{ Cell = 'D45'; Value = '0.6815' }
    @{ Cell = 'E45'; Value = '  -2.97%  ' }
    @{ Cell = 'D46'; Value = '2.387' }
    @{ Cell = 'E46'; Value = '  -4.70%  ' }
    @{ Cell = 'D47'; Value = '4.055' }
    @{ Cell = 'E47'; Value = '  -1.05%  ' }
    @{ Cell = 'D48'; Value = '0.9999' }
    @{ Cell = 'E48'; Value = '  +0.18%  ' }
    @{ Cell = 'D49'; Value = '0.08173' }
    @{ Cell = 'E49'; Value = '  -1.58%  ' }
    @{ Cell = 'D50'; Value = '133.50' }
    @{ Cell = 'E50'; Value = '  -1.29%  ' }
    @{ Cell = 'D51'; Value = '1.207' }
    @{ Cell = 'E51'; Value = '  -2.57%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.ClearFormats()
}
